$d = $word.ActiveDocument

# --- 1) Collapse "Ergebnisueberpruefung..." paragraph into two clean runs ---
$p2 = $d.Paragraphs(2)
$p2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Ergebnisüberprüfung pro Klasse scikit image Überlagerung des Originalbildes und des Zahlenplans</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> Heatmap</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($p2Xml)

# --- 2) Insert three new detail paragraphs right after it ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs(3)
$newPara1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:ind w:left="1440"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Reshaping array to 11 arrays of size (990, 820, 4) RGBA with values per pixel (0, 0, 1, probability)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara1.Range.InsertXML($newPara1Xml)

$newPara1 = $d.Paragraphs(3)
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(4)
$newPara2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:ind w:left="1440"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Creating image by using PIL.Image.fromarray</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara2.Range.InsertXML($newPara2Xml)

$newPara2 = $d.Paragraphs(4)
$newPara2.Range.InsertParagraphAfter()
$newPara3 = $d.Paragraphs(5)
$newPara3Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:ind w:left="1440"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Overlaying with original plan by using image.paste</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara3.Range.InsertXML($newPara3Xml)

# --- 3) Collapse "Rechtecke einzeichnen (...)" paragraph into a single run ---
$pRechtecke = $d.Paragraphs(8)
$pRechteckeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Rechtecke einzeichnen (opencv Polygon, opencv kleinstmögliches Rechteck) und die Ergebnisse in Textform abspeichern</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pRechteckeText = $pRechtecke.Range.Text
if ($pRechteckeText -notlike "Rechtecke einzeichnen*") {
    throw "Unexpected paragraph at index 8: $pRechteckeText"
}
$pRechtecke.Range.InsertXML($pRechteckeXml)

# --- 4) Collapse "Fully-Convolutional-Network einbauen" paragraph into a single run ---
$pFcn = $d.Paragraphs(10)
$pFcnXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Fully-Convolutional-Network einbauen</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pFcnText = $pFcn.Range.Text
if ($pFcnText -notlike "Fully-Convolutional*") {
    throw "Unexpected paragraph at index 10: $pFcnText"
}
$pFcn.Range.InsertXML($pFcnXml)

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
